$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 492-493, pushing the existing data (old rows
# 492..622) down to 494..624. This also extends the used range / dimension
# from R622 to R624, matching the target diff.
$ws.Range("A492:A493").EntireRow.Insert()

# Populate the newly inserted rows with the new weekly observation
# (Coliflor, Terminal La Palmera de La Serena) for date serial 44642.

# Row 492 - Calidad "Primera"
$ws.Cells.Item(492, 1).Value = 8
$ws.Cells.Item(492, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(492, 3).Value = "Coquimbo"
$ws.Cells.Item(492, 4).Value = 44642
$ws.Cells.Item(492, 5).Value = 4
$ws.Cells.Item(492, 6).Value = 100112008
$ws.Cells.Item(492, 7).Value = "Coliflor"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 2200
$ws.Cells.Item(492, 11).Value = 950
$ws.Cells.Item(492, 12).Value = 1000
$ws.Cells.Item(492, 13).Value = 975
$ws.Cells.Item(492, 14).Value = "`$/unidad"
$ws.Cells.Item(492, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(492, 16).Value = 975
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"

# Row 493 - Calidad "Segunda"
$ws.Cells.Item(493, 1).Value = 8
$ws.Cells.Item(493, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(493, 3).Value = "Coquimbo"
$ws.Cells.Item(493, 4).Value = 44642
$ws.Cells.Item(493, 5).Value = 4
$ws.Cells.Item(493, 6).Value = 100112008
$ws.Cells.Item(493, 7).Value = "Coliflor"
$ws.Cells.Item(493, 8).Value = "Sin especificar"
$ws.Cells.Item(493, 9).Value = "Segunda"
$ws.Cells.Item(493, 10).Value = 1200
$ws.Cells.Item(493, 11).Value = 850
$ws.Cells.Item(493, 12).Value = 900
$ws.Cells.Item(493, 13).Value = 875
$ws.Cells.Item(493, 14).Value = "`$/unidad"
$ws.Cells.Item(493, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(493, 16).Value = 875
$ws.Cells.Item(493, 17).Value = 1
$ws.Cells.Item(493, 18).Value = "Hortaliza"
